# "A02SZL1_Fin Buff Calc" update — refreshed input figures from the 502
# Part C / Part L source report.
#
#   D3 (Enter Gross Expenditures From 502 Part C): 547650.06 -> 548401.97
#   D5 (Enter Total Labor Cost From 502 Part L):   297809.10 -> 298007.16
#
# Everything else on the sheet (D6 tech-labor %, D8/E8 added cost, D9/E9
# suggested standard added cost) is formula-driven off D3/D5 and
# recalculates automatically. The active selection also moved from D4 to
# D5, matching where the author left the cursor after editing the labor
# cost figure.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 548401.97
$ws.Range("D5").Value = 298007.15999999997

$ws.Range("D5").Select()

$wb.Save()
